$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "29.868.18"
Set-TextValue "E2" "  +0.95%  "
Set-TextValue "D3" "1.624.68"
Set-TextValue "E3" "  +1.16%  "
Set-TextValue "E4" "  -0.36%  "
Set-TextValue "D5" "214.24"
Set-TextValue "E5" "  +0.85%  "
Set-TextValue "D6" "0.521"
Set-TextValue "E6" "  -0.45%  "
Set-TextValue "E7" "  -0.39%  "
Set-TextValue "D8" "29.64"
Set-TextValue "E8" "  +10.59%  "
Set-TextValue "E9" "  +3.21%  "
Set-TextValue "E10" "  +1.62%  "
Set-TextValue "E11" "  +0.46%  "
Set-TextValue "D12" "1.856.75"
Set-TextValue "E12" "  +1.14%  "
Set-TextValue "D13" "1.619.89"
Set-TextValue "E13" "  +0.63%  "
Set-TextValue "E14" "  +5.97%  "
Set-TextValue "D16" "29.913.86"
Set-TextValue "D17" "8.83"
Set-TextValue "E17" "  +16.25%  "
Set-TextValue "D18" "64.66"
Set-TextValue "E18" "  +1.93%  "
Set-TextValue "D19" "244.42"
Set-TextValue "E19" "  +1.65%  "
Set-TextValue "D20" "0.0₃0705"
Set-TextValue "E20" "  +1.79%  "
Set-TextValue "E21" "  -0.31%  "
Set-TextValue "E22" "  +3.50%  "
Set-TextValue "D23" "9.63"
Set-TextValue "E23" "  +4.51%  "
Set-TextValue "E24" "  +2.62%  "
Set-TextValue "E25" "  +1.45%  "
Set-TextValue "E26" "  +2.60%  "
Set-TextValue "E27" "  +1.80%  "
Set-TextValue "E28" "  +3.17%  "
Set-TextValue "D29" "0.996"
Set-TextValue "E29" "  -0.34%  "
Set-TextValue "E30" "  +3.40%  "
Set-TextValue "E31" "  +5.28%  "
Set-TextValue "E32" "  +3.57%  "
Set-TextValue "D33" "3.23"
Set-TextValue "E33" "  +3.73%  "
Set-TextValue "D34" "1.425.57"
Set-TextValue "E34" "  +1.26%  "
Set-TextValue "E35" "  +6.88%  "
Set-TextValue "E36" "  -0.08%  "
Set-TextValue "E37" "  +2.13%  "
Set-TextValue "E38" "  -0.72%  "
Set-TextValue "D40" "0.557"
Set-TextValue "E40" "  +3.54%  "
Set-TextValue "D41" "0.0508"
Set-TextValue "E41" "  +2.85%  "
Set-TextValue "E42" "  +0.34%  "
Set-TextValue "E43" "  +4.82%  "
Set-TextValue "D44" "54.15"
Set-TextValue "E44" "  +0.07%  "
Set-TextValue "D45" "69.29"
Set-TextValue "E46" "  +18.45%  "
Set-TextValue "E48" "  +2.78%  "
Set-TextValue "D49" "1.765.14"
Set-TextValue "E49" "  +1.08%  "
Set-TextValue "D50" "88.38"
Set-TextValue "E50" "  +1.95%  "
Set-TextValue "D51" "0.0₆0107"
Set-TextValue "E51" "  +3.54%  "
